# Update info capture from master_url file
# Applies: insert "Type" column (D) and "Short Name" column (K), populate new
# homologation data, rename headers, drop the old hard-coded hyperlinks, and
# move the active selection to A2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Insert the two new columns, carrying neighbouring formatting with them.
# ---------------------------------------------------------------------------

# New "Type" column, pushed in front of the old "Linea" column (old D).
$ws.Columns("D:D").Insert()

# New "Short Name" column, pushed in front of the "Link" column (now K).
$ws.Columns("K:K").Insert()

# ---------------------------------------------------------------------------
# 2. Header row text updates.
# ---------------------------------------------------------------------------
$ws.Range("B1").Value = "Homologo Mansfield"
$ws.Range("D1").Value = "Type"
$ws.Range("J1").Value = "Descripcion"
$ws.Range("K1").Value = "Short Name"

# ---------------------------------------------------------------------------
# 3. Fill in the new "Type" / "Linea" / "Rough in" / "Short Name" data.
# ---------------------------------------------------------------------------

# Row 2 - Mansfield VX1 ADA EL Bowl
$ws.Range("D2").Value = "Bowl"
$ws.Range("E2").Value = "VX1"
$ws.Range("F2").Value = "'12"
$ws.Range("K2").Value = "Mansfield VX1 ADA EL Bowl"

# Row 3 - Mansfield VX1 1.28 gpf Tank
$ws.Range("D3").Value = "Tank"
$ws.Range("E3").Value = "VX1"
$ws.Range("F3").Value = "NA"
$ws.Range("G3").Value = "NA"
$ws.Range("K3").Value = "Mansfield VX1 1,28 gpf Tank"

# Row 4 - Mansfield Summit 1.28 gpf Tank
$ws.Range("D4").Value = "Tank"
$ws.Range("E4").Value = "Summit"
$ws.Range("F4").Value = "NA"
$ws.Range("G4").Value = "NA"
$ws.Range("K4").Value = "Mansfield Summit 1,28 gpf Tank"

# ---------------------------------------------------------------------------
# 4. Re-apply header formatting: the whole header row is now center-aligned
#    (it used to be a mix of left/center), and the new D1 cell needs the
#    "no-left-border" header look that the Linea/Rough-in group uses.
# ---------------------------------------------------------------------------
$ws.Range("A1:L1").HorizontalAlignment = -4108  # xlCenter

$ws.Range("E1").Copy() | Out-Null
$ws.Range("D1").PasteSpecial(-4122)             # xlPasteFormats

# ---------------------------------------------------------------------------
# 5. Data-row formatting: columns A-I are center aligned, the two
#    description columns (Descripcion / Short Name) are left aligned like
#    the old "Name" column was, and the Link column keeps the hyperlink look.
# ---------------------------------------------------------------------------
$ws.Range("D2:D4").HorizontalAlignment = -4108  # xlCenter
$ws.Range("E2:E4").HorizontalAlignment = -4108
$ws.Range("F2:F4").HorizontalAlignment = -4108
$ws.Range("I2:I4").HorizontalAlignment = -4108

$ws.Range("J2:J4").Copy() | Out-Null
$ws.Range("K2:K4").PasteSpecial(-4122)          # xlPasteFormats

$ws.Application.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 6. Drop the old explicit hyperlinks - the Link column is now plain text.
# ---------------------------------------------------------------------------
$ws.Hyperlinks.Delete()

# ---------------------------------------------------------------------------
# 7. Column widths, to roughly match the post-edit auto-fit layout.
# ---------------------------------------------------------------------------
$ws.Columns("D:D").ColumnWidth = 13.7109375
$ws.Columns("E:E").ColumnWidth = 14.42578125
$ws.Columns("F:F").ColumnWidth = 16.28515625
$ws.Columns("G:G").ColumnWidth = 16.42578125
$ws.Columns("H:H").ColumnWidth = 20.7109375
$ws.Columns("I:I").ColumnWidth = 20.7109375
$ws.Columns("J:J").ColumnWidth = 61.28515625
$ws.Columns("K:K").ColumnWidth = 29.42578125
$ws.Columns("L:L").ColumnWidth = 100.28515625

# ---------------------------------------------------------------------------
# 8. Move the active selection to A2, matching the saved view state.
# ---------------------------------------------------------------------------
$ws.Range("A2").Select()
